# Apply cryptos list price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.713.03'
$ws.Range('E2').Value = '  +1.21%  '

# Row 3
$ws.Range('D3').Value = '1.875.14'
$ws.Range('E3').Value = '  +1.54%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '332.04'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +3.60%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.004'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.04%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4719'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +6.52%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3956'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.50%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.57'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.38%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08033'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.23%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.023'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.11%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.84'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.05%  '

# Row 13
$ws.Range('D13').Value = '1.892.80'
$ws.Range('E13').Value = '  +2.37%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.957'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.37%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.148'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.20%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.007'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.02%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001048'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.88%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '87.19'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.27%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06648'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.49%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.28'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.35%  '

# Row 21
$ws.Range('E21').Value = '  -0.11%  '

# Row 22
$ws.Range('D22').Value = '27.790.71'
$ws.Range('E22').Value = '  +1.55%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.503'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.11%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.02'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.89%  '

# Row 25
$ws.Range('E25').Value = '  +1.61%  '

# Row 26
$ws.Range('D26').Value = '2.117.69'
$ws.Range('E26').Value = '  +2.49%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '156.18'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.17%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.23'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.30%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.100'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +3.92%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.579'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.68%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '122.52'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.00%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9676'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.10%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09551'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.16%  '

# Row 34
$ws.Range('E34').Value = '  -0.63%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.632'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.01%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.303'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +2.37%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06109'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +3.29%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02261'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.58%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.230'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.88%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.187'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.29%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.003'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.02%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.5987'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.16%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1910'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.94%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '10.25'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.93%  '

# Row 45
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5706'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.78%  '

# Row 46
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.248'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.93%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '12.26'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +1.63%  '

# Row 48
$ws.Range('E48').Value = '  +1.72%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.932'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.35%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06818'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.25%  '

# Row 51
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.00000000312'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +9.02%  '
